$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add a new row of data: youtube platform with a new URL (plain text, not a hyperlink)
$ws.Range("A5").Value = "youtube"
$ws.Range("B5").Value = "https://www.youtube.com/watch?v=whCINxitNkM"

# Match the recorded selection after the edit
$ws.Range("A12").Select()
